# TC for Login added
#
# - Rename the "Login" sheet to "Login&Logout" (test cases for Logout were
#   added to that sheet).
# - "Login&Logout" becomes the active/selected sheet tab (it was the
#   "User Registration" sheet before).
# - Selection on "User Registration" moves to C15:C17.
# - Selection on "Login&Logout" moves to B22, scrolled so row 6 is at the
#   top of the view.

$wb = $excel.ActiveWorkbook

$loginSheet   = $wb.Worksheets.Item("Login")
$userRegSheet = $wb.Worksheets.Item("User Registration")

# Rename "Login" -> "Login&Logout"
$loginSheet.Name = "Login&Logout"

# "User Registration" is no longer the active tab; update its lingering
# selection to C15:C17.
$userRegSheet.Activate()
$userRegSheet.Range("C15:C17").Select()

# Make "Login&Logout" the active tab and set its view/selection.
$loginSheet.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$loginSheet.Range("B22").Select()
